# Update odds values in row 2 of the active sheet to reflect the latest
# FlashScore data refresh (Jogos da Semana FlashScore 2024-11-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.29
$ws.Range("H2").Value = 5.3
$ws.Range("J2").Value = 1.72
$ws.Range("K2").Value = 2.62
$ws.Range("L2").Value = 7

$ws.Range("P2").Value = 4.85
$ws.Range("Q2").Value = 1.45
$ws.Range("R2").Value = 2.55
$ws.Range("S2").Value = 1.26
$ws.Range("T2").Value = 3.45
$ws.Range("U2").Value = 1.8
$ws.Range("V2").Value = 1.91
$ws.Range("W2").Value = 9
$ws.Range("X2").Value = 7.3
$ws.Range("Y2").Value = 8.75
$ws.Range("Z2").Value = 8.5
$ws.Range("AA2").Value = 10.25
$ws.Range("AB2").Value = 23
$ws.Range("AD2").Value = 10.75

$ws.Range("AG2").Value = 450
$ws.Range("AH2").Value = 28
$ws.Range("AJ2").Value = 26
$ws.Range("AL2").Value = 90
$ws.Range("AM2").Value = 65

$ws.Range("AS2").Value = 175
$ws.Range("AT2").Value = 3.45

$ws.Range("AW2").Value = 9.5
$ws.Range("AX2").Value = 45
$ws.Range("AY2").Value = 40
$ws.Range("AZ2").Value = 300
$ws.Range("BA2").Value = 250
$ws.Range("BB2").Value = 450
